# "upload mapa y fake"
# Populates the (previously empty) "Sheet3" worksheet with a per-department
# COVID case table ("mapa"/"fake" data export) and updates the active
# sheet/selection bookkeeping so Sheet3 becomes the tab shown on open.

$wb = $excel.ActiveWorkbook
$wsData = $wb.Worksheets.Item(1)   # "Sheet2" - existing daily summary sheet
$wsMapa = $wb.Worksheets.Item(2)   # "Sheet3" - empty sheet to populate

# --- Header row ---------------------------------------------------------
$wsMapa.Range("A1").Value = "Departamento"
$wsMapa.Range("B1").Value = "Positivos"
$wsMapa.Range("C1").Value = "Fecha"

# --- Data rows: Departamento / Positivos / Fecha ------------------------
$rows = @(
    ,("LIMA", 58, 43905)
    ,("CALLAO", 2, 43905)
    ,("ANCASH", 1, 43905)
    ,("AREQUIPA", 2, 43905)
    ,("CUSCO", 1, 43905)
    ,("HUANUCO", 2, 43905)
    ,("ICA", 1, 43905)
    ,("LA LIBERTAD", 1, 43905)
    ,("LAMBAYEQUE", 1, 43905)
    ,("PIURA", 2, 43905)
    ,("LIMA", 70, 43906)
    ,("CALLAO", 3, 43906)
    ,("LAMBAYEQUE", 3, 43906)
    ,("AREQUIPA", 2, 43906)
    ,("HUANUCO", 2, 43906)
    ,("PIURA", 2, 43906)
    ,("ANCASH", 1, 43906)
    ,("CUSCO", 1, 43906)
    ,("ICA", 1, 43906)
    ,("LA LIBERTAD", 1, 43906)
    ,("LIMA", 37, 43904)
    ,("AREQUIPA", 2, 43904)
    ,("CUSCO", 1, 43904)
    ,("HUANUCO", 2, 43904)
    ,("ICA", 1, 43904)
    ,("LIMA", 32, 43903)
    ,("AREQUIPA", 2, 43903)
    ,("CUSCO", 1, 43903)
    ,("HUANUCO", 2, 43903)
    ,("ICA", 1, 43903)
    ,("LIMA", 22, 43902)
    ,("LIMA", 17, 43901)
    ,("LIMA", 11, 43900)
    ,("LIMA", 9, 43899)
    ,("LIMA", 7, 43898)
    ,("LIMA", 1, 43897)
    ,("LIMA", 6, 43896)
)

$r = 2
foreach ($row in $rows) {
    $dept = $row[0]
    $pos = $row[1]
    $fecha = $row[2]

    $wsMapa.Cells.Item($r, 1).Value = $dept
    $wsMapa.Cells.Item($r, 2).Value = $pos
    $wsMapa.Cells.Item($r, 3).Value = $fecha
    $wsMapa.Cells.Item($r, 3).NumberFormat = "d-mmm"

    $r = $r + 1
}

# --- Column sizing --------------------------------------------------------
$wsMapa.Columns.Item(1).ColumnWidth = 12.1666666

# --- Selection / active-tab bookkeeping -----------------------------------
# Sheet2 keeps its own remembered selection but is no longer the active tab.
$wsData.Activate() | Out-Null
$wsData.Range("E2:E3").Select() | Out-Null

# Sheet3 becomes the active/selected tab, with B41 as the remembered selection.
$wsMapa.Activate() | Out-Null
$wsMapa.Range("B41").Select() | Out-Null
